# Onderzoek-data.xlsx - "wijzigingen david doorgevoerd"
#
# The C1:C4 "address" header labels (City/Street/House_Number/Postal, shown
# in red) are removed from column C and re-entered (with corrected, lower-
# case / snake_case names) as new rows inserted into column A, right after
# the "isadmin" row. That pushes the rest of column A (created_at ... amount)
# down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing column-A entries (rows 11-20: created_at ... amount)
#    down by four rows to make room for the four new rows (11-14).
#    Walk bottom-up so we never clobber a value before it's been copied.
for ($r = 20; $r -ge 11; $r--) {
    $srcVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 4, 1).Value = $srcVal
}

# 2) Remove the old red-font labels from column C (C1:C4) - contents only,
#    the column formatting / style stays in place.
$ws.Range("C1").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# 3) Write the corrected labels into the freshly-opened rows 11-14 of
#    column A, using a font distinct from the old red one (plain/automatic
#    colour rather than red).
$ws.Range("A11").Value = "city"
$ws.Range("A12").Value = "street"
$ws.Range("A13").Value = "house_nr"
$ws.Range("A14").Value = "postalcode"

$newLabels = $ws.Range("A11:A14")
$newLabels.Font.Color = 0

# 4) Update the sheet's selection to match the author's cursor position when
#    they saved the file.
$ws.Range("D14").Select() | Out-Null
